# Applies the 2023-12-06 GitHub Actions crypto-price refresh to Sheet1.
# Column D ("Price") holds numeric-looking text (e.g. "44.200.27", "0.635") that must
# stay plain text, so those cells are pre-formatted as Text ("@") before the assignment
# and the explicit style is cleared back to Normal afterwards so no stray formatting sticks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '44.200.27'
$ws.Range('E2').Value = '  +5.73%  '

# Row 3
$ws.Range('D3').Value = '2.266.77'
$ws.Range('E3').Value = '  +2.57%  '

# Row 4
$ws.Range('E4').Value = '  -0.17%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '230.79'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.17%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.635'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +2.85%  '

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '63.76'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +5.18%  '

# Row 8
$ws.Range('E8').Value = '  -0.02%  '

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.440'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +9.68%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.102'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +14.31%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '56.39'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -1.27%  '

# Row 12: Avalanche -> TRON (rank swap)
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.106'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +2.27%  '

# Row 13: TRON -> Avalanche (rank swap)
$ws.Range('B13').Value = 'Avalanche'
$ws.Range('C13').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '25.75'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +16.47%  '

# Row 14
$ws.Range('D14').Value = '2.597.19'
$ws.Range('E14').Value = '  +2.36%  '

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '15.71'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +2.00%  '

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '5.98'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +7.40%  '

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.826'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +3.76%  '

# Row 18
$ws.Range('D18').Value = '2.259.51'
$ws.Range('E18').Value = '  +2.47%  '

# Row 19
$ws.Range('D19').Value = '43.976.76'
$ws.Range('E19').Value = '  +5.01%  '

# Row 20
$ws.Range('E20').Value = '  +10.48%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '73.53'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +2.22%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.03'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -0.42%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '255.30'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +5.35%  '

# Row 25
$ws.Range('E25').Value = '  +3.39%  '

# Row 26
$ws.Range('E26').Value = '  -1.22%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.00'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +4.28%  '

# Row 28: Monero -> WEMIXToken (rank swap)
$ws.Range('B28').Value = 'WEMIXToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '3.25'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +22.99%  '

# Row 29: WEMIXToken -> Monero (rank swap)
$ws.Range('B29').Value = 'Monero'
$ws.Range('C29').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '172.04'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +2.18%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '20.85'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +5.64%  '

# Row 31
$ws.Range('E31').Value = '  -1.33%  '

# Row 32
$ws.Range('E32').Value = '  -2.87%  '

# Row 33
$ws.Range('E33').Value = '  +2.72%  '

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.0683'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +5.54%  '

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '4.73'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +3.08%  '

# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '4.88'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -1.44%  '

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '3.85'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +8.60%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '6.74'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +6.34%  '

# Row 39
$ws.Range('E39').Value = '  -0.18%  '

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.0257'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +6.71%  '

# Row 41
$ws.Range('E41').Value = '  -0.01%  '

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '8.36'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -2.12%  '

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '17.42'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +8.26%  '

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0962'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +1.11%  '

# Row 45: Aave -> FTXToken (rank swap)
$ws.Range('B45').Value = 'FTXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '4.42'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +3.09%  '

# Row 46: FTXToken -> Aave (rank swap)
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '97.54'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +0.97%  '

# Row 47
$ws.Range('E47').Value = '  -0.25%  '

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.000210'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -11.63%  '

# Row 49
$ws.Range('D49').Value = '1.448.09'
$ws.Range('E49').Value = '  -0.29%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.28'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +2.43%  '

# Row 51
$ws.Range('E51').Value = '  +1.34%  '
